$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("teleostei post-hatching")

# Rename header labels: swap "Method Of Euthanasia" to BV1 stays text-same,
# but BT1/BU1 get renamed from "Sampling Data ..." to "Sampling Day ..."
$ws.Range("BT1").Value = "Sampling Day Start Time"
$ws.Range("BU1").Value = "Sampling Day End Time"
$ws.Range("BV1").Value = "Method Of Euthanasia"

# Scroll the view so column BG is the top-left visible column and select BV1
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = $ws.Range("BG1").Column
$ws.Range("BV1").Select()
